# Fixed exclude path bug in path_list file
$wb = $excel.ActiveWorkbook

# --- "include" sheet: drop the extra path row (row 3), keep row 2 as-is ---
$wsInclude = $wb.Worksheets.Item("include")
$wsInclude.Rows.Item(3).Delete()
$wsInclude.Range("A3:XFD3").Select()

# --- "exclude" sheet: correct the path value on row 2, drop the now-stale rows 3 & 4 ---
$wsExclude = $wb.Worksheets.Item("exclude")
$wsExclude.Range("B2").Value = "C:\Temp\images1\folder2"
$wsExclude.Rows.Item(4).Delete()
$wsExclude.Rows.Item(3).Delete()
$wsExclude.Range("C7").Select()
